$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Senior Golang Developer", "https://www.dice.com/job-detail/501e2625-4e56-4074-81e3-7e181ce193a4", "Hybrid in Sunnyvale, California", "Contract", "Depends on Experience", "P2PSoftTek Inc"),
    @("Golang Developer (W2)", "https://www.dice.com/job-detail/2fd7b3dd-4e05-4030-bbe5-876fbf976cee", "Hybrid in Sunnyvale, California", "Contract", "$40 - $65", "Lorvenk Technologies LLC"),
    @("Senior .Net/GoLang Integration Developer - Primarily Remote", "https://www.dice.com/job-detail/8228c51c-e3b3-42ad-8344-a098563b08de", "Remote", "Contract", "Depends on Experience", "EmployVision")
)

$startRow = 27
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
